$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Resistor footprint/supplier swap (rows 7 & 8) ---
# Row 7: was Vishay Dale CRCW2010100RFKEFHP (SMD 2010 100R) -> Stackpole RSMF2JT330R (through hole 330R, 2W)
$ws.Range("A7").Value = "Stackpole Electronics RSMF2JT330R "
$ws.Range("B7").Value = "Through Hole 330 ohm 2 watt resistor"
$ws.Range("G7").Value = "https://www.digikey.com/product-detail/en/stackpole-electronics-inc/RSF2JT330R/RSF2JT330RCT-ND/2021796"

# Row 8: was Vishay Dale CRCW201010R0FKEF (SMD 2010 10R) -> Stackpole CF12JT10R0 (through hole 10R, 1/2W)
$ws.Range("A8").Value = "Stackpole Electronics  CF12JT10R0 "
$ws.Range("B8").Value = "Through Hole 10 ohm ½ watt Resistor"
$ws.Range("G8").Value = "https://www.digikey.com/product-detail/en/stackpole-electronics-inc/CF12JT10R0/CF12JT10R0CT-ND/1830446"

# --- Updated unit prices for rows 7 & 8 ---
$ws.Range("C7").Value = 0.29
$ws.Range("C8").Value = 0.1

# --- D5 formula now rounds the quantity up to a whole unit ---
$ws.Range("D5").Formula = "=_xlfn.CEILING.MATH(73/20)"

# --- Selection moved to A10 ---
$ws.Range("A10").Select()
